$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Anexo 3 report: add an extra data column/value.
# Header row swaps meaning: A1 -> "codigo_completo", B1 -> "nombre"
$ws.Range("A1").Value = "codigo_completo"
$ws.Range("B1").Value = "nombre"

# The longer header text now needs a wider column (best-fit sizing).
$ws.Columns.Item(1).ColumnWidth = 15.666666666666666

# Move the active selection to A2 (below the header row).
$ws.Range("A2").Select()
